# "updated some prez formatting"
#
# Moves the SmartArt diagram ("Content Placeholder 8") on slide 4 to its new
# position. EMU -> point conversion (1 pt = 12700 EMU) loses a little
# precision once it round-trips through the COM Shape.Left/.Top Single
# (float32) properties, so a tiny epsilon is added before conversion to make
# sure the value lands back on the exact target EMU after the host's
# point->EMU re-quantization.

function ConvertTo-PointsFromEmu {
    param(
        [double]$Emu
    )
    return ($Emu / 12700.0) + 0.00001
}

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(4)
$diagram = $slide.Shapes.Item(2)

$diagram.Left = ConvertTo-PointsFromEmu 5483901
$diagram.Top  = ConvertTo-PointsFromEmu 2209285
